$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "41.219.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.435.59"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.60%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "88.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.26%  "

$ws.Range("E7").Value = "  -2.12%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  -4.14%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "32.21"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.45%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0833"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.32%  "

$ws.Range("E12").Value = "  -2.85%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.805.93"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.86%  "

$ws.Range("E14").Value = "  -2.92%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.70"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.08%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.427.72"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.774"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.07%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "41.148.99"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.16%  "

$ws.Range("E19").Value = "  -3.55%  "

$ws.Range("E20").Value = "  -3.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "71.98"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.72%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.97%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.22%  "

$ws.Range("E25").Value = "  -0.01%  "

$ws.Range("E26").Value = "  -2.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.59%  "

$ws.Range("E28").Value = "  -3.24%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.57"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.24%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "34.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "157.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.12%  "

$ws.Range("E32").Value = "  +0.01%  "

$ws.Range("E33").Value = "  -4.87%  "

$ws.Range("E34").Value = "  -2.22%  "

$ws.Range("E35").Value = "  -3.20%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.93"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.34%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "16.64"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -5.03%  "

$ws.Range("E38").Value = "  -0.59%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.53%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.78%  "

$ws.Range("E42").Value = "  -7.02%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.983.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.02%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0276"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.53%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.37"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.31%  "

$ws.Range("E46").Value = "  -5.23%  "

$ws.Range("E47").Value = "  +2.58%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.664.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "95.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.29%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "73.23"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.11%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "51.69"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.83%  "

